$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.279.72'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.791.64'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.02'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.29'
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0690'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.048.12'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.01'
$ws.Range("E13").Value = '  -4.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.789.05'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.226.56'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.03'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0803'
$ws.Range("E19").Value = '  +2.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '246.74'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.42'
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0521'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.77'
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.84'
$ws.Range("E33").Value = '  +5.57%  '
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.441.56'
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.61'
$ws.Range("E36").Value = '  +8.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.662'
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0189'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.51'
$ws.Range("E40").Value = '  +1.89%  '
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.07'
$ws.Range("E42").Value = '  +5.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.73'
$ws.Range("E43").Value = '  +1.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.924'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("E47").Value = '  +0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.942.43'
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.57'
$ws.Range("E49").Value = '  -2.35%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("E51").Value = '  -7.05%  '
